$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New glucose readings appended after the existing log (rows 348-401),
# continuing the 5-minute interval series from 2026/02/12 19:03 through 23:28.
$rows = @(
    @(348, "2026/02/12 19:03", "9.4"),
    @(349, "2026/02/12 19:08", "9.7"),
    @(350, "2026/02/12 19:13", "9.5"),
    @(351, "2026/02/12 19:18", "9.4"),
    @(352, "2026/02/12 19:23", "8.5"),
    @(353, "2026/02/12 19:28", "7.8"),
    @(354, "2026/02/12 19:33", "8.0"),
    @(355, "2026/02/12 19:38", "7.9"),
    @(356, "2026/02/12 19:43", "8.0"),
    @(357, "2026/02/12 19:48", "8.3"),
    @(358, "2026/02/12 19:53", "8.4"),
    @(359, "2026/02/12 19:58", "8.1"),
    @(360, "2026/02/12 20:03", "8.3"),
    @(361, "2026/02/12 20:08", "9.7"),
    @(362, "2026/02/12 20:13", "11.4"),
    @(363, "2026/02/12 20:18", "13.4"),
    @(364, "2026/02/12 20:23", "15.2"),
    @(365, "2026/02/12 20:28", "16.3"),
    @(366, "2026/02/12 20:33", "16.1"),
    @(367, "2026/02/12 20:38", "17.4"),
    @(368, "2026/02/12 20:43", "17.1"),
    @(369, "2026/02/12 20:48", "18.1"),
    @(370, "2026/02/12 20:53", "18.3"),
    @(371, "2026/02/12 20:58", "19.3"),
    @(372, "2026/02/12 21:03", "20.1"),
    @(373, "2026/02/12 21:08", "20.6"),
    @(374, "2026/02/12 21:13", "21.6"),
    @(375, "2026/02/12 21:18", "21.3"),
    @(376, "2026/02/12 21:23", "22.4"),
    @(377, "2026/02/12 21:28", "22.1"),
    @(378, "2026/02/12 21:33", "22.4"),
    @(379, "2026/02/12 21:38", "22.1"),
    @(380, "2026/02/12 21:43", "22.5"),
    @(381, "2026/02/12 21:48", "22.1"),
    @(382, "2026/02/12 21:53", "22.0"),
    @(383, "2026/02/12 21:58", "22.1"),
    @(384, "2026/02/12 22:03", "22.5"),
    @(385, "2026/02/12 22:08", "21.8"),
    @(386, "2026/02/12 22:13", "22.3"),
    @(387, "2026/02/12 22:18", "22.1"),
    @(388, "2026/02/12 22:23", "21.0"),
    @(389, "2026/02/12 22:28", "20.4"),
    @(390, "2026/02/12 22:33", "19.8"),
    @(391, "2026/02/12 22:38", "19.4"),
    @(392, "2026/02/12 22:43", "19.1"),
    @(393, "2026/02/12 22:48", "18.5"),
    @(394, "2026/02/12 22:53", "18.6"),
    @(395, "2026/02/12 22:58", "19.0"),
    @(396, "2026/02/12 23:03", "19.3"),
    @(397, "2026/02/12 23:08", "19.5"),
    @(398, "2026/02/12 23:13", "19.5"),
    @(399, "2026/02/12 23:18", "19.7"),
    @(400, "2026/02/12 23:23", "20.3"),
    @(401, "2026/02/12 23:28", "20.0")
)

foreach ($row in $rows) {
    $r = $row[0]
    $timeVal = $row[1]
    $glucoseVal = $row[2]
    $ws.Range("A$r").Value = $timeVal
    # Leading apostrophe forces text storage so numeric-looking readings
    # (e.g. '9.4') stay shared-string text instead of becoming Number cells.
    $ws.Range("B$r").Value = "'$glucoseVal"
}

# Strip the transient quote-prefix formatting picked up above so the column
# keeps the workbook default (General) style, same as every other data cell.
$ws.Range("B348:B401").ClearFormats()
